$d = $word.ActiveDocument

$pairs = @(
    @("2025-10-12 Sunday", "2025-10-13 Monday"),
    @("28+28=56", "86-47=39"),
    @("8+25=33", "32-16=16"),
    @("14+77=91", "12+74=86"),
    @("83-33=50", "27+20=47"),
    @("12+30=42", "3+90=93"),
    @("42-42=0", "49+9=58"),
    @("78+9=87", "86-27=59"),
    @("57+28=85", "27+49=76"),
    @("2+32=34", "71-2=69"),
    @("38-3=35", "35+22=57"),
    @("30+64=94", "91-69=22"),
    @("13+58=71", "38-28=10"),
    @("55+29=84", "86-80=6"),
    @("90-19=71", "12+64=76"),
    @("71+3=74", "56+8=64"),
    @("18-10=8", "65+10=75"),
    @("87-33=54", "65+17=82"),
    @("2+63=65", "83+11=94"),
    @("43+26=69", "34+55=89"),
    @("32+43=75", "16+28=44"),
    @("2+56=58", "98+1=99"),
    @("25-6=19", "5+58=63"),
    @("44-27=17", "61-34=27"),
    @("13-11=2", "32+56=88"),
    @("88-41=47", "3+29=32"),
    @("22-8=14", "73-41=32"),
    @("58+37=95", "78-56=22"),
    @("67-31=36", "18-4=14"),
    @("12+61=73", "37+44=81"),
    @("67+19=86", "73+12=85"),
    @("75-71=4", "32-3=29"),
    @("69-2=67", "86-22=64"),
    @("42-2=40", "82-6=76"),
    @("49+10=59", "64-33=31"),
    @("87-65=22", "97-12=85"),
    @("18-17=1", "99-19=80"),
    @("94-62=32", "87-82=5"),
    @("16+27=43", "85-27=58"),
    @("30+4=34", "52+40=92"),
    @("99-78=21", "90-44=46"),
    @("76-57=19", "25+1=26"),
    @("23+64=87", "28+24=52"),
    @("75+11=86", "81-43=38"),
    @("8+90=98", "70-55=15"),
    @("43+39=82", "46-16=30"),
    @("37+19=56", "99-95=4"),
    @("29-16=13", "98-54=44"),
    @("11+61=72", "91-58=33"),
    @("56-43=13", "41+45=86"),
    @("16-12=4", "71+24=95"),
    @("3+8=11", "43+17=60"),
    @("24+19=43", "8+58=66"),
    @("15+69=84", "31+31=62"),
    @("83-22=61", "25+3=28"),
    @("48-35=13", "89-56=33"),
    @("51-46=5", "25+15=40"),
    @("54-15=39", "36+25=61"),
    @("49-29=20", "12+39=51"),
    @("44+18=62", "62+4=66"),
    @("91-31=60", "48-25=23"),
    @("86+12=98", "24-1=23"),
    @("20-10=10", "36-30=6"),
    @("28-15=13", "63+5=68"),
    @("45+45=90", "71-15=56"),
    @("27+21=48", "58-2=56"),
    @("27+30=57", "78+0=78"),
    @("16-0=16", "63-41=22"),
    @("55+34=89", "89-42=47"),
    @("4-2=2", "41-19=22"),
    @("97-34=63", "95-7=88"),
    @("12+72=84", "51+27=78"),
    @("34+57=91", "72-15=57"),
    @("87-4=83", "20+3=23"),
    @("45+25=70", "23+58=81"),
    @("23-15=8", "12+7=19"),
    @("61-55=6", "44+33=77"),
    @("36+59=95", "63-17=46"),
    @("18+53=71", "53+44=97"),
    @("43-3=40", "83-25=58"),
    @("51-47=4", "65+1=66"),
    @("29+57=86", "44-37=7"),
    @("83-72=11", "42-29=13"),
    @("18+78=96", "8+68=76"),
    @("2+7=9", "36+62=98"),
    @("3+63=66", "8+53=61"),
    @("34+50=84", "16+20=36"),
    @("37+46=83", "0+37=37"),
    @("80+5=85", "92-84=8"),
    @("44+31=75", "40-13=27"),
    @("18+14=32", "96-0=96"),
    @("84-61=23", "80-49=31"),
    @("73-39=34", "39+10=49"),
    @("34+44=78", "50-19=31"),
    @("81-0=81", "99-38=61"),
    @("63+1=64", "11+71=82"),
    @("80-59=21", "19+12=31"),
    @("34+62=96", "24+39=63"),
    @("4+84=88", "98-30=68"),
    @("59-21=38", "55+3=58"),
    @("39+29=68", "18+38=56"),
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
    }
}

Write-Output "Done"
